$d = $word.ActiveDocument

$newText = "El restaurante ha estado empeorando en su gestión de pedidos debido a la falta de organización y control que se observa en el proceso de toma de pedidos. Esta situación ha generado inconvenientes tanto para el personal encargado de esta área, como para los clientes, quienes han tenido que esperar más tiempo por sus pedidos o incluso cancelarlos por la demora. Esta mala experiencia ha afectado la reputación y la fidelidad del restaurante, lo que pone en riesgo su sostenibilidad y competitividad en el mercado."

$p2 = $d.Paragraphs(2)
$p3 = $d.Paragraphs(3)

# Replace the text spanning both paragraphs 2 and 3 with the new merged text.
$rng = $d.Range($p2.Range.Start, $p3.Range.End)
$rng.Text = $newText

# Remove the now-empty trailing paragraph (formerly paragraph 3), merging
# paragraph 2's pilcrow with paragraph 3's so only one paragraph remains.
$p2b = $d.Paragraphs(2)
$p3b = $d.Paragraphs(3)
$delRng = $d.Range($p2b.Range.End - 1, $p3b.Range.End)
$delRng.Delete()
